$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 4 oldest months (2020-10 .. 2021-01) from the top of the data block
$ws.Rows("2:5").Delete() | Out-Null

# Delete the 4 trailing rows left over at the bottom (previously rows 52:55, now rows 48:51)
$ws.Rows("48:51").Delete() | Out-Null

# Rewrite every remaining data row (2..51) with the refreshed metrics
$data = @(
  ,@("2021-02", 25, 94, 26.59574468085106)
  ,@("2021-03", 30, 77, 38.96103896103897)
  ,@("2021-04", 37, 93, 39.78494623655914)
  ,@("2021-05", 43, 88, 48.86363636363637)
  ,@("2021-06", 52, 91, 57.14285714285714)
  ,@("2021-07", 63, 122, 51.63934426229508)
  ,@("2021-08", 61, 100, 61.0)
  ,@("2021-09", 63, 121, 52.06611570247934)
  ,@("2021-10", 67, 122, 54.91803278688525)
  ,@("2021-11", 77, 134, 57.46268656716418)
  ,@("2021-12", 78, 127, 61.41732283464567)
  ,@("2022-01", 79, 117, 67.52136752136752)
  ,@("2022-02", 85, 143, 59.44055944055944)
  ,@("2022-03", 91, 143, 63.63636363636363)
  ,@("2022-04", 89, 146, 60.95890410958904)
  ,@("2022-05", 77, 139, 55.39568345323741)
  ,@("2022-06", 88, 112, 78.57142857142857)
  ,@("2022-07", 81, 117, 69.23076923076923)
  ,@("2022-08", 71, 104, 68.26923076923077)
  ,@("2022-09", 76, 108, 70.37037037037037)
  ,@("2022-10", 69, 127, 54.33070866141733)
  ,@("2022-11", 74, 98, 75.51020408163265)
  ,@("2022-12", 75, 123, 60.97560975609756)
  ,@("2023-01", 75, 104, 72.11538461538461)
  ,@("2023-02", 78, 124, 62.90322580645162)
  ,@("2023-03", 80, 121, 66.11570247933885)
  ,@("2023-04", 92, 130, 70.76923076923077)
  ,@("2023-05", 107, 141, 75.88652482269504)
  ,@("2023-06", 103, 143, 72.02797202797203)
  ,@("2023-07", 101, 138, 73.18840579710145)
  ,@("2023-08", 108, 142, 76.05633802816901)
  ,@("2023-09", 115, 160, 71.875)
  ,@("2023-10", 115, 156, 73.71794871794873)
  ,@("2023-11", 110, 164, 67.07317073170732)
  ,@("2023-12", 99, 150, 66.0)
  ,@("2024-01", 105, 143, 73.42657342657343)
  ,@("2024-02", 123, 182, 67.58241758241759)
  ,@("2024-03", 121, 161, 75.15527950310559)
  ,@("2024-04", 136, 189, 71.95767195767195)
  ,@("2024-05", 143, 197, 72.58883248730965)
  ,@("2024-06", 141, 209, 67.46411483253588)
  ,@("2024-07", 150, 215, 69.76744186046511)
  ,@("2024-08", 165, 225, 73.33333333333333)
  ,@("2024-09", 172, 226, 76.10619469026548)
  ,@("2024-10", 169, 227, 74.44933920704845)
  ,@("2024-11", 151, 242, 62.39669421487604)
  ,@("2024-12", 123, 197, 62.43654822335025)
  ,@("2025-01", 145, 195, 74.35897435897436)
  ,@("2025-02", 160, 216, 74.07407407407408)
  ,@("2025-03", 86, 230, 37.39130434782609)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
}

Write-Output ("Final rows: " + $ws.UsedRange.Rows.Count)
